# Fix typo in syllabus: "147-158" -> "147-171"
#
# The target text lives in a single run ("147-158"). Word's actual
# revision-tracked edit flow (replace the trailing "58" with "71")
# ends up producing two runs sharing the same run properties
# ("147-1" and "71") once the insertion is accepted, instead of
# silently re-merging the whole paragraph's runs. So: turn on track
# changes, perform the edit as a delete + insert, accept the
# resulting revisions, then restore the original track-changes state.

$d = $word.ActiveDocument

$originalTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $true

$target = $d.Content
$target.Find.Execute("147-158", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($target.Find.Found) {
    $end = $target.End

    # Remove the erroneous trailing "58", leaving "147-1".
    $tail = $d.Range($end - 2, $end)
    $tail.Delete()

    # Insert the correction "71" right after, as its own run.
    $insertionPoint = $d.Range($end - 2, $end - 2)
    $insertionPoint.InsertAfter("71")

    # Accept the tracked insert/delete so the document ends up clean.
    for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
        $d.Revisions.Item($i).Accept()
    }
}

$d.TrackRevisions = $originalTrackRevisions
